# Refresh the cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
# Price/volume cells are stored as plain text in this sheet, so any replacement value
# that looks like a bare number is written with a leading apostrophe - Excel's normal
# "force text" convention - so it round-trips as text instead of being coerced to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.154.69'
$ws.Range('E2').Value = '  +1.68%  '
# Row 3
$ws.Range('D3').Value = '2.589.51'
$ws.Range('E3').Value = '  +0.37%  '
# Row 4
$ws.Range('E4').Value = '  -0.03%  '
# Row 5
$ws.Range('D5').Value = '''566.53'
$ws.Range('E5').Value = '  +0.64%  '
# Row 6
$ws.Range('D6').Value = '''141.56'
$ws.Range('E6').Value = '  -0.63%  '
# Row 7
$ws.Range('E7').Value = '  -0.19%  '
# Row 8
$ws.Range('E8').Value = '  -0.08%  '
# Row 9
$ws.Range('D9').Value = '2.609.16'
$ws.Range('E9').Value = '  +0.85%  '
# Row 10
$ws.Range('D10').Value = '''6.57'
$ws.Range('E10').Value = '  -1.02%  '
# Row 11
$ws.Range('E11').Value = '  +1.29%  '
# Row 12
$ws.Range('E12').Value = '  +6.20%  '
# Row 13
$ws.Range('E13').Value = '  -6.06%  '
# Row 14
$ws.Range('D14').Value = '3.051.47'
$ws.Range('E14').Value = '  +0.64%  '
# Row 15
$ws.Range('D15').Value = '60.150.97'
$ws.Range('E15').Value = '  +1.72%  '
# Row 16
$ws.Range('E16').Value = '  +1.23%  '
# Row 17
$ws.Range('D17').Value = '''0.0000139'
$ws.Range('E17').Value = '  +2.00%  '
# Row 18
$ws.Range('D18').Value = '2.600.16'
$ws.Range('E18').Value = '  +0.62%  '
# Row 19
$ws.Range('D19').Value = '''11.33'
$ws.Range('E19').Value = '  +9.32%  '
# Row 20
$ws.Range('E20').Value = '  +1.87%  '
# Row 21
$ws.Range('D21').Value = '''345.45'
$ws.Range('E21').Value = '  +2.61%  '
# Row 22
$ws.Range('E22').Value = '  +8.37%  '
# Row 23
$ws.Range('E23').Value = '  -0.01%  '
# Row 24
$ws.Range('E24').Value = '  +15.87%  '
# Row 25
$ws.Range('D25').Value = '''62.95'
$ws.Range('E25').Value = '  -1.86%  '
# Row 26
$ws.Range('E26').Value = '  -0.40%  '
# Row 27
$ws.Range('E27').Value = '  -2.02%  '
# Row 28
$ws.Range('D28').Value = '''7.65'
$ws.Range('E28').Value = '  +4.74%  '
# Row 29
$ws.Range('E29').Value = '  +0.85%  '
# Row 30
$ws.Range('E30').Value = '  +7.32%  '
# Row 31
$ws.Range('E31').Value = '  -0.07%  '
# Row 32
$ws.Range('D32').Value = '''6.31'
$ws.Range('E32').Value = '  +3.67%  '
# Row 33
$ws.Range('D33').Value = '''160.55'
$ws.Range('E33').Value = '  -0.57%  '
# Row 34
$ws.Range('D34').Value = '''19.41'
$ws.Range('E34').Value = '  +2.59%  '
# Row 35
$ws.Range('D35').Value = '''4.22'
$ws.Range('E35').Value = '  +5.14%  '
# Row 36
$ws.Range('D36').Value = '''0.956'
$ws.Range('E36').Value = '  +9.42%  '
# Row 37
$ws.Range('E37').Value = '  +4.55%  '
# Row 38
$ws.Range('E38').Value = '  +7.53%  '
# Row 39
$ws.Range('D39').Value = '''37.71'
$ws.Range('E39').Value = '  +0.71%  '
# Row 40
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D40').Value = '''0.851'
$ws.Range('E40').Value = '  -2.62%  '
# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '''3.81'
$ws.Range('E41').Value = '  +4.18%  '
# Row 42
$ws.Range('D42').Value = '''293.65'
$ws.Range('E42').Value = '  -0.20%  '
# Row 43
$ws.Range('D43').Value = '''138.01'
$ws.Range('E43').Value = '  +4.82%  '
# Row 44
$ws.Range('E44').Value = '  -0.26%  '
# Row 45
$ws.Range('E45').Value = '  +1.44%  '
# Row 46
$ws.Range('E46').Value = '  +0.54%  '
# Row 47
$ws.Range('D47').Value = '''19.55'
$ws.Range('E47').Value = '  +2.85%  '
# Row 48
$ws.Range('E48').Value = '  +1.58%  '
# Row 49
$ws.Range('E49').Value = '  +8.81%  '
# Row 50
$ws.Range('E50').Value = '  +2.72%  '
# Row 51
$ws.Range('E51').Value = '  +0.10%  '
